$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill row 24 of Tableau1 with a new journal entry
# (set end-time / start-time before the date so the dependent "Durée"
# formula picks up both operands before it recalculates)
$ws.Range("G24").Value = 0.625
$ws.Range("F24").Value = 0.57638888888888895
$ws.Range("E24").Value = 44265

$ws.Range("I24").Value = "Développement"
$ws.Range("L24").Value = "Faire que le programme vérifie les coordonnée et quil l'affiche"
$ws.Range("J24").Value = "Faire une vérification des coordonnées et l'affichage"
$ws.Range("K24").Value = "CPNV"

$excel.CalculateFull()

# Apply the date / time number formats used by the other rows of the table
# (copy just the formatting, so the existing cell styles are reused)
$ws.Range("E23:G23").Copy()
$ws.Range("E24:G24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row grows to fit the wrapped "Tâche"/"Descriptif" text, same as the row above it
$ws.Rows.Item(24).RowHeight = $ws.Rows.Item(23).RowHeight

# Update the selection to match the saved view state
$ws.Range("J25").Select()
